{"js": "// Replace each known old value with its corresponding new value.\n// The document is a \"two-digit multiplied by two-digit\" worksheet: a date\n// heading paragraph plus a table of \"A\u00d7B=C\" cells. Every old string in the\n// mapping is unique in the document, so body.search(...) + insertText(...,\n// replace) on the single match is sufficient and preserves the run's\n// existing formatting (font/size) because insertText replaces only the\n// text content of the matched range.\nconst replacements = [\n  [\"2024-03-15 Friday\", \"2024-03-16 Saturday\"],\n  [\"89\u00d728=2492\", \"14\u00d757=798\"],\n  [\"66\u00d720=1320\", \"73\u00d757=4161\"],\n  [\"56\u00d772=4032\", \"31\u00d762=1922\"],\n  [\"76\u00d777=5852\", \"80\u00d720=1600\"],\n  [\"16\u00d791=1456\", \"77\u00d760=4620\"],\n  [\"68\u00d717=1156\", \"85\u00d796=8160\"],\n  [\"32\u00d713=416\", \"16\u00d724=384\"],\n  [\"75\u00d714=1050\", \"92\u00d717=1564\"],\n  [\"90\u00d727=2430\", \"76\u00d714=1064\"],\n  [\"93\u00d716=1488\", \"87\u00d789=7743\"],\n  [\"52\u00d779=4108\", \"14\u00d762=868\"],\n  [\"21\u00d717=357\", \"43\u00d799=4257\"],\n  [\"53\u00d785=4505\", \"65\u00d767=4355\"],\n  [\"54\u00d721=1134\", \"93\u00d788=8184\"],\n  [\"14\u00d724=336\", \"65\u00d757=3705\"],\n  [\"21\u00d715=315\", \"93\u00d730=2790\"],\n  [\"58\u00d735=2030\", \"40\u00d714=560\"],\n  [\"91\u00d739=3549\", \"65\u00d747=3055\"],\n  [\"76\u00d729=2204\", \"70\u00d739=2730\"],\n  [\"88\u00d752=4576\", \"90\u00d763=5670\"],\n  [\"46\u00d745=2070\", \"14\u00d714=196\"],\n  [\"76\u00d713=988\", \"57\u00d787=4959\"],\n  [\"88\u00d769=6072\", \"63\u00d749=3087\"],\n  [\"28\u00d792=2576\", \"35\u00d737=1295\"],\n  [\"34\u00d799=3366\", \"73\u00d770=5110\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Replace every match (expected to be exactly one per mapping entry).\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each known old value (\"A\u00d7B=C\" cell text, plus the date heading)\n# with its corresponding new value, using Word's Find/Replace (wdReplaceOne)\n# on $d.Content so only the exact matched text is rewritten and the run's\n# existing formatting (font/size) stays untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = '2024-03-15 Friday'; New = '2024-03-16 Saturday' },\n    @{ Old = '89\u00d728=2492'; New = '14\u00d757=798' },\n    @{ Old = '66\u00d720=1320'; New = '73\u00d757=4161' },\n    @{ Old = '56\u00d772=4032'; New = '31\u00d762=1922' },\n    @{ Old = '76\u00d777=5852'; New = '80\u00d720=1600' },\n    @{ Old = '16\u00d791=1456'; New = '77\u00d760=4620' },\n    @{ Old = '68\u00d717=1156'; New = '85\u00d796=8160' },\n    @{ Old = '32\u00d713=416'; New = '16\u00d724=384' },\n    @{ Old = '75\u00d714=1050'; New = '92\u00d717=1564' },\n    @{ Old = '90\u00d727=2430'; New = '76\u00d714=1064' },\n    @{ Old = '93\u00d716=1488'; New = '87\u00d789=7743' },\n    @{ Old = '52\u00d779=4108'; New = '14\u00d762=868' },\n    @{ Old = '21\u00d717=357'; New = '43\u00d799=4257' },\n    @{ Old = '53\u00d785=4505'; New = '65\u00d767=4355' },\n    @{ Old = '54\u00d721=1134'; New = '93\u00d788=8184' },\n    @{ Old = '14\u00d724=336'; New = '65\u00d757=3705' },\n    @{ Old = '21\u00d715=315'; New = '93\u00d730=2790' },\n    @{ Old = '58\u00d735=2030'; New = '40\u00d714=560' },\n    @{ Old = '91\u00d739=3549'; New = '65\u00d747=3055' },\n    @{ Old = '76\u00d729=2204'; New = '70\u00d739=2730' },\n    @{ Old = '88\u00d752=4576'; New = '90\u00d763=5670' },\n    @{ Old = '46\u00d745=2070'; New = '14\u00d714=196' },\n    @{ Old = '76\u00d713=988'; New = '57\u00d787=4959' },\n    @{ Old = '88\u00d769=6072'; New = '63\u00d749=3087' },\n    @{ Old = '28\u00d792=2576'; New = '35\u00d737=1295' },\n    @{ Old = '34\u00d799=3366'; New = '73\u00d770=5110' }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $($pair.Old)\"\n    }\n}\n\nWrite-Output \"done\"\n"}
